$wb = $excel.ActiveWorkbook

# --- Hoja3: add "Mayor" column (D) with MAX(C2:C5), and repoint A2's
#     denominator from $C5 to the new $D2 ---
$ws3 = $wb.Worksheets.Item("Hoja3")

$ws3.Range("D1").Value = "Mayor"
$ws3.Range("D2").Formula = "=MAX(C2:C5)"
$ws3.Range("A2").Formula = "=(SUM(`$B2:`$B5)/`$D2)"

# --- Active sheet moves from Hoja5 to Hoja3 ---
$ws3.Activate()
